$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text formatting so numeric-looking strings
# (e.g. "309.84") are not auto-converted into real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.064.40'
$ws.Range("E2").Value = '  -3.11%  '
$ws.Range("D3").Value = '1.715.09'
$ws.Range("E3").Value = '  -3.02%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '309.84'
$ws.Range("E5").Value = '  -5.74%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4785'
$ws.Range("E7").Value = '  +4.92%  '
$ws.Range("D8").Value = '0.3449'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '42.18'
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").Value = '0.07255'
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").Value = '  -5.07%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '19.77'
$ws.Range("E13").Value = '  -4.62%  '
$ws.Range("D14").Value = '5.832'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").Value = '1.714.27'
$ws.Range("E15").Value = '  -3.04%  '
$ws.Range("D16").Value = '6.825'
$ws.Range("E16").Value = '  -5.04%  '
$ws.Range("D17").Value = '87.29'
$ws.Range("E17").Value = '  -5.73%  '
$ws.Range("D18").Value = '0.00001034'
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").Value = '0.06377'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '16.46'
$ws.Range("E21").Value = '  -3.01%  '
$ws.Range("D22").Value = '5.609'
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("D23").Value = '27.118.59'
$ws.Range("E23").Value = '  -2.99%  '
$ws.Range("D24").Value = '10.74'
$ws.Range("E24").Value = '  -4.23%  '
$ws.Range("D25").Value = '2.098'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '20.02'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("D27").Value = '150.82'
$ws.Range("E27").Value = '  -5.26%  '
$ws.Range("D28").Value = '1.908.22'
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("D29").Value = '2.056'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("D30").Value = '120.69'
$ws.Range("E30").Value = '  -2.67%  '
$ws.Range("D31").Value = '1.029'
$ws.Range("E31").Value = '  -4.99%  '
$ws.Range("D32").Value = '0.09248'
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("D33").Value = '3.600'
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").Value = '5.319'
$ws.Range("E34").Value = '  -5.38%  '
$ws.Range("D35").Value = '1.467'
$ws.Range("E35").Value = '  +5.75%  '
$ws.Range("D36").Value = '0.02183'
$ws.Range("E36").Value = '  -4.18%  '
$ws.Range("D37").Value = '0.05850'
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("D38").Value = '10.90'
$ws.Range("E38").Value = '  -7.95%  '
$ws.Range("D39").Value = '0.1988'
$ws.Range("E39").Value = '  -4.90%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '4.715'
$ws.Range("E41").Value = '  -4.77%  '
$ws.Range("D42").Value = '0.5956'
$ws.Range("E42").Value = '  -4.82%  '
$ws.Range("D43").Value = '1.085'
$ws.Range("E43").Value = '  -8.01%  '
$ws.Range("D44").Value = '7.509'
$ws.Range("E44").Value = '  -3.72%  '
$ws.Range("D45").Value = '12.72'
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").Value = '3.579'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").Value = '0.5574'
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("D48").Value = '118.79'
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").Value = '1.823'
$ws.Range("E49").Value = '  -5.64%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.099'
$ws.Range("E50").Value = '  -3.06%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.06631'
$ws.Range("E51").Value = '  -2.94%  '
